$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the formatting of the last existing row (row 7) down into the
# two new rows so the new cells inherit the same date/time number formats
# (styles) already used by the table, then insert them as real rows.
$ws.Rows("7").Copy()
$ws.Rows("8").Insert(-4121)

$ws.Rows("7").Copy()
$ws.Rows("9").Insert(-4121)

# New data: 11/10/2013 (1:00) and 12/10/2013 (4:30)
$ws.Range("A8").Value = 41559
$ws.Range("B8").Value = 0.041666666666666664

$ws.Range("A9").Value = 41560
$ws.Range("B9").Value = 0.1875

# Update the visible selection to match the extended range.
[void]$ws.Range("B4:B9").Select()
